# Tambah kolom "deadline" (tanggal) setelah "tahun" dan kolom "konversi"
# sebelum "denda", sesuai permintaan revisi tim kampus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sisipkan kolom baru "deadline" tepat sebelum kolom "nim" (kolom D).
#    Ini otomatis menggeser semua kolom sesudahnya (nim..denda) satu kolom
#    ke kanan dan menyesuaikan referensi formula yang ada.
$ws.Columns("D").Insert()

# 2) Sisipkan kolom baru "konversi" tepat sebelum kolom terakhir ("denda",
#    yang sekarang berada di kolom Q setelah penyisipan di atas).
$ws.Columns("Q").Insert()

# Isi kolom "konversi" lebih dulu supaya string bersama (shared string)
# "konversi" dialokasikan sebelum "deadline", sama seperti urutan penulis asli.
$ws.Range("Q1").Value = "konversi"
$ws.Range("Q2").Value = 999999

# Isi header dan nilai kolom "deadline" dengan format tanggal.
$ws.Range("D1").Value = "deadline"
$deadlineCell = $ws.Range("D2")
$deadlineCell.Value2 = 44196
$deadlineCell.NumberFormat = "mm-dd-yy"

# Lebarkan kolom "deadline" otomatis mengikuti isinya.
$ws.Columns("D").AutoFit()

# Tulis ulang rumus SUM yang dibagi (shared formula) pada baris 3:4 supaya
# jangkauannya mencakup kolom baru (G3:G4 -> SUM(H:R) alih-alih SUM(G:P)).
$ws.Range("G3:G4").Formula = "=SUM(H3:R3)"

# Pindahkan sel aktif seperti pada perubahan asli.
$null = $ws.Range("D3").Select()
